# Apply updated Fitness (column C) values on Sheet1, as described by the diff.
# Each entry is (startRow, endRow, newValue); the same value is written to every
# row in the inclusive range on column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ranges = @(
    @(2, 2, 12428),
    @(3, 3, 12312),
    @(4, 7, 10057),
    @(8, 10, 9419),
    @(11, 17, 9098),
    @(18, 23, 8853),
    @(24, 25, 8743),
    @(26, 27, 8581),
    @(28, 40, 8264),
    @(41, 42, 8188),
    @(43, 57, 8048),
    @(58, 75, 7937),
    @(169, 177, 7651)
)

foreach ($r in $ranges) {
    $startRow = $r[0]
    $endRow = $r[1]
    $value = $r[2]
    $rangeAddress = "C$startRow`:C$endRow"
    $ws.Range($rangeAddress).Value = $value
}
